$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# Add a new "version" column (D) to the settings sheet, with value 1,
# mirroring the other XLSForm settings columns (form_title, form_id,
# allow_choice_duplicates).
$ws.Range("D1").Value = "version"
$ws.Range("D2").Value = 1

# Move the selection the way a user would after typing the new values
# (down into D2, then on to the next empty row).
$null = $ws.Range("D3").Select()
